# Applies scheduled-runner updates to Hades_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Generated from the authoritative OOXML diff; each Range.Value assignment updates a single cell
# to match the post-edit numeric value. Cells removed in the diff (no longer present in the row)
# are cleared by assigning $null so the saved XML omits the <c> element, matching the target state.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = $null
$ws.Range("H137").Value = 2566528.2
$ws.Range("I137").Value = 5001662
$ws.Range("J137").Value = 3229.5789
$ws.Range("K137").Value = 15004986
$ws.Range("L137").Value = 9688.736699999999
$ws.Range("M137").Value = -15002436
$ws.Range("N137").Value = -14788.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2037.91
$ws.Range("I32").Value = 1828
$ws.Range("J32").Value = 4826.7144
$ws.Range("K32").Value = 1828
$ws.Range("L32").Value = 4826.7144
$ws.Range("M32").Value = -1541
$ws.Range("N32").Value = -5400.7144
$ws.Range("H61").Value = 24440334
$ws.Range("I61").Value = 31282382
$ws.Range("K61").Value = 31282382
$ws.Range("M61").Value = -31282170
$ws.Range("H74").Value = 5729329
$ws.Range("I74").Value = 8656418
$ws.Range("J74").Value = 70290
$ws.Range("K74").Value = 8656418
$ws.Range("L74").Value = 70290
$ws.Range("M74").Value = -8655544
$ws.Range("N74").Value = -72038
$ws.Range("H77").Value = 5729329
$ws.Range("I77").Value = 8656418
$ws.Range("J77").Value = 70290
$ws.Range("K77").Value = 43282090
$ws.Range("L77").Value = 351450
$ws.Range("M77").Value = -43277722
$ws.Range("N77").Value = -360186
$ws.Range("H122").Value = 3705379.2
$ws.Range("I122").Value = 1686.6538
$ws.Range("K122").Value = 5059.9614
$ws.Range("M122").Value = -2609.9614
$ws.Range("H125").Value = 59965
$ws.Range("J125").Value = 59965
$ws.Range("L125").Value = 59965
$ws.Range("N125").Value = -69805
$ws.Range("H132").Value = 36146.793
$ws.Range("I132").Value = 26218.6
$ws.Range("J132").Value = 58209.445
$ws.Range("K132").Value = 78655.79999999999
$ws.Range("L132").Value = 174628.335
$ws.Range("M132").Value = -76125.79999999999
$ws.Range("N132").Value = -179688.335
$ws.Range("H136").Value = 24440334
$ws.Range("I136").Value = 31282382
$ws.Range("K136").Value = 93847146
$ws.Range("M136").Value = -93844596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2331.0833
$ws.Range("I134").Value = 1724.3103
$ws.Range("J134").Value = 4844.857
$ws.Range("K134").Value = 5172.9309
$ws.Range("L134").Value = 14534.571
$ws.Range("M134").Value = -2637.9309
$ws.Range("N134").Value = -19604.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1985.4333
$ws.Range("I31").Value = 1407.7826
$ws.Range("K31").Value = 1407.7826
$ws.Range("M31").Value = -1112.7826
$ws.Range("H34").Value = 1985.4333
$ws.Range("I34").Value = 1407.7826
$ws.Range("K34").Value = 1407.7826
$ws.Range("M34").Value = -1205.7826
$ws.Range("H58").Value = 18520066
$ws.Range("I58").Value = 27028100
$ws.Range("J58").Value = 2583.353
$ws.Range("K58").Value = 27028100
$ws.Range("L58").Value = 2583.353
$ws.Range("M58").Value = -27027897
$ws.Range("N58").Value = -2989.353
$ws.Range("H107").Value = 369.05554
$ws.Range("I107").Value = 338.75
$ws.Range("J107").Value = 429.66666
$ws.Range("K107").Value = 338.75
$ws.Range("L107").Value = 429.66666
$ws.Range("M107").Value = 1581.25
$ws.Range("N107").Value = -4269.66666
$ws.Range("H132").Value = 32354.787
$ws.Range("I132").Value = 1431.2727
$ws.Range("J132").Value = 94201.82000000001
$ws.Range("K132").Value = 4293.8181
$ws.Range("L132").Value = 282605.46
$ws.Range("M132").Value = -1763.8181
$ws.Range("N132").Value = -287665.46
$ws.Range("H134").Value = 24092.36
$ws.Range("I134").Value = 1544.619
$ws.Range("K134").Value = 4633.857
$ws.Range("M134").Value = -2098.857
$ws.Range("H136").Value = 18520066
$ws.Range("I136").Value = 27028100
$ws.Range("J136").Value = 2583.353
$ws.Range("K136").Value = 81084300
$ws.Range("L136").Value = 7750.059
$ws.Range("M136").Value = -81081750
$ws.Range("N136").Value = -12850.059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 37037596
$ws.Range("I12").Value = 111112350
$ws.Range("J12").Value = 215.05556
$ws.Range("K12").Value = 333337050
$ws.Range("L12").Value = 645.16668
$ws.Range("M12").Value = -333336877
$ws.Range("N12").Value = -991.16668
$ws.Range("H113").Value = 502.89655
$ws.Range("I113").Value = 445.29166
$ws.Range("J113").Value = 779.4
$ws.Range("K113").Value = 1335.87498
$ws.Range("L113").Value = 2338.2
$ws.Range("M113").Value = 834.1250199999999
$ws.Range("N113").Value = -6678.2
$ws.Range("H131").Value = 24578.816
$ws.Range("I131").Value = 554.5454999999999
$ws.Range("J131").Value = 34366.48
$ws.Range("K131").Value = 1663.6365
$ws.Range("L131").Value = 103099.44
$ws.Range("M131").Value = 3376.3635
$ws.Range("N131").Value = -113179.44

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 62170.03
$ws.Range("I132").Value = 48514.76
$ws.Range("J132").Value = 86066.75
$ws.Range("K132").Value = 145544.28
$ws.Range("L132").Value = 258200.25
$ws.Range("M132").Value = -143014.28
$ws.Range("N132").Value = -263260.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1730.2
$ws.Range("I40").Value = 1681.8125
$ws.Range("J40").Value = 1923.75
$ws.Range("K40").Value = 1681.8125
$ws.Range("L40").Value = 1923.75
$ws.Range("M40").Value = -1545.8125
$ws.Range("N40").Value = -2195.75
$ws.Range("H100").Value = 49192.57
$ws.Range("I100").Value = 112205.555
$ws.Range("J100").Value = 1932.8334
$ws.Range("K100").Value = 112205.555
$ws.Range("L100").Value = 1932.8334
$ws.Range("M100").Value = -111664.555
$ws.Range("N100").Value = -3014.8334
$ws.Range("H122").Value = 2747.7
$ws.Range("I122").Value = 2789.6
$ws.Range("J122").Value = 2622
$ws.Range("K122").Value = 8368.799999999999
$ws.Range("L122").Value = 7866
$ws.Range("M122").Value = -5918.799999999999
$ws.Range("N122").Value = -12766
$ws.Range("H132").Value = 20957
$ws.Range("I132").Value = 1434.275
$ws.Range("J132").Value = 86032.75
$ws.Range("K132").Value = 4302.825000000001
$ws.Range("L132").Value = 258098.25
$ws.Range("M132").Value = -1772.825000000001
$ws.Range("N132").Value = -263158.25
$ws.Range("H136").Value = 63644.375
$ws.Range("I136").Value = 40944.64
$ws.Range("J136").Value = 144714.86
$ws.Range("K136").Value = 122833.92
$ws.Range("L136").Value = 434144.58
$ws.Range("M136").Value = -120283.92
$ws.Range("N136").Value = -439244.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 47602.117
$ws.Range("I132").Value = 30982.787
$ws.Range("J132").Value = 102445.9
$ws.Range("K132").Value = 92948.361
$ws.Range("L132").Value = 307337.7
$ws.Range("M132").Value = -90418.361
$ws.Range("N132").Value = -312397.7
$ws.Range("H136").Value = 40298.668
$ws.Range("I136").Value = 23693.908
$ws.Range("J136").Value = 144671.42
$ws.Range("K136").Value = 71081.724
$ws.Range("L136").Value = 434014.26
$ws.Range("M136").Value = -68531.724
$ws.Range("N136").Value = -439114.26
